$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 / Spring 2022 block (rows 4-10) ---
# Re-order the Spring 2022 course list and insert a new Spring 2022 row,
# then append a new "CPSC 4000" course (0 credits) as its own row.
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("C5").Value = "CPSC 4135"

$ws.Range("A6").Value = "CPSC 3121"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CPSC 4148"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "CPSC 4000"
$ws.Range("B7").Value = 0

# --- Fall 2023 / Spring 2023 block (rows 13-19) ---
# Re-order the course list: swap Fall 2023 entries, move one Spring 2023
# entry into its own Fall 2023 row, and move CPSC 4175 under Fall 2023 too.
$ws.Range("A13").Value = "CPSC 4155"
$ws.Range("C13").Value = "CPSC 4176"

$ws.Range("A14").Value = "CPSC 4157"
$ws.Range("C14:D14").Value = ""

$ws.Range("A15").Value = "CPSC 4175"
$ws.Range("B15").Value = 3

# --- Fall 2024 / Spring 2024 / Summer 2024 block (rows 22-24) ---
# These courses are no longer listed under Fall 2024 (moved above), so clear them.
$ws.Range("A22:F24").Value = ""

# --- Remove the Fall 2025 / Spring 2025 / Summer 2025 and
#     Fall 2026 / Spring 2026 / Summer 2026 blocks entirely (rows 30-47) ---
$ws.Rows("30:47").Delete()
